$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Mayo de 2020 a las 12:05"

# Row data: [row, Ciudad, Casos totales, Casos activos, Recuperados, Muertes]
$rows = @(
    @(4, 'Madrid', 66005, 40383, 16813, 8809),
    @(5, 'Cataluña', 55464, 25904, 23678, 5882),
    @(6, 'Castilla y Leon', 18272, 7672, 8670, 1930),
    @(7, 'Castilla-La Mancha', 16513, 6294, 7349, 2870),
    @(8, 'Pais Vasco', 13156, 14646, 0, 1418),
    @(9, 'Andalucia', 12401, 10148, 909, 1344),
    @(10, 'Bizkaia/Vizcaya', 10332, 7124, 4423, 815),
    @(11, 'Galicia', 9323, 8283, 439, 601),
    @(12, 'Valencia/Valencia', 7011, 4066, 2767, 658),
    @(13, 'Ciudad Real', 6414, 1996, 3370, 1048),
    @(14, 'Aragon', 5432, 3534, 1061, 837),
    @(15, 'Zaragoza', 5125, 2618, 1854, 653),
    @(16, 'Navarra', 5116, 3550, 1066, 500),
    @(17, 'Araba/Alava', 4868, 7124, 4423, 356),
    @(18, 'Alacant/Alicante', 4671, 3017, 1938, 470),
    @(19, 'Valladolid', 4322, 1538, 2430, 354),
    @(20, 'Salamanca', 4124, 1167, 2603, 354),
    @(21, 'La Rioja', 4016, 2927, 741, 348),
    @(22, 'Malaga', 3919, 2098, 1548, 273),
    @(23, 'Toledo', 3831, 1929, 1156, 746),
    @(24, 'Albacete', 3756, 1363, 1888, 505),
    @(25, 'Leon', 3526, 1575, 1550, 401),
    @(26, 'Segovia', 3398, 857, 2341, 200),
    @(27, 'Gipuzkoa/Guipuzcoa', 3116, 7124, 4423, 283),
    @(28, 'Sevilla', 3062, 1711, 1078, 273),
    @(29, 'Granada', 2995, 2412, 308, 275),
    @(30, 'Extremadura', 2919, 2422, 10, 487),
    @(31, 'Burgos', 2721, 891, 1625, 205),
    @(32, 'Asturias', 2359, 1053, 996, 310),
    @(33, 'Gran Canaria', 2280, 1506, 623, 151),
    @(34, 'Soria', 2276, 390, 1767, 119),
    @(35, 'Cantabria', 2246, 1981, 62, 203),
    @(36, 'Castello/Castellon', 1988, 1098, 699, 200),
    @(37, 'Caceres', 1970, 1411, 159, 400),
    @(38, 'A Coruña', 1969, 333, 1788, 67),
    @(39, 'Avila', 1900, 613, 1155, 132),
    @(40, 'Jaen', 1724, 1103, 451, 170),
    @(41, 'Cordoba', 1656, 1308, 243, 105),
    @(42, 'Pontevedra', 1536, 333, 1411, 30),
    @(43, 'Tenerife', 1532, 897, 529, 106),
    @(44, 'Murcia', 1508, 1782, 0, 139),
    @(45, 'Cadiz', 1477, 511, 824, 142),
    @(46, 'Guadalajara', 1255, 365, 641, 249),
    @(47, 'Cuenca', 1214, 591, 319, 304),
    @(48, 'Palencia', 1185, 324, 781, 80),
    @(49, 'Huesca', 1061, 439, 523, 99),
    @(50, 'Badajoz', 955, 1071, 0, 92),
    @(51, 'Zamora', 897, 317, 496, 84),
    @(52, 'Ourense', 751, 333, 660, 22),
    @(53, 'Almeria', 677, 447, 180, 50),
    @(54, 'Teruel', 643, 363, 197, 83),
    @(55, 'Lugo', 586, 333, 520, 11),
    @(56, 'Huelva', 510, 328, 134, 48),
    @(57, 'Mallorca', 210, 18, 194, 12),
    @(58, 'Ceuta', 125, 98, 23, 4),
    @(59, 'Melilla', 119, 115, 2, 2),
    @(60, 'La Palma', 80, 60, 17, 3),
    @(61, 'Lanzarote', 74, 60, 9, 5),
    @(62, 'Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
    @(63, 'Fuerteventura', 23, 22, 1, 0),
    @(64, 'Ibiza', 21, 18, 20, 1),
    @(65, 'Menorca', 15, 18, 13, 0),
    @(66, 'Arroyo de la Luz', 7, 0, 7, 0),
    @(67, 'La Gomera', 7, 7, 0, 0),
    @(68, 'El Hierro', 1, 1, 0, 0),
    @(69, 'Formentera', 0, 10, 0, 8)

)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
